$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7906
$ws.Range("L2").Value = 5342
$ws.Range("L3").Value = 5762
$ws.Range("L4").Value = 1410
$ws.Range("L5").Value = 343
$ws.Range("L6").Value = 4783
$ws.Range("K7").Value = 27584
$ws.Range("L7").Value = 17640

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L6").Value = 134
$ws.Range("L7").Value = 574
$ws.Range("L11").Value = 291
$ws.Range("L15").Value = 136
$ws.Range("L19").Value = 479
$ws.Range("L20").Value = 436
$ws.Range("L24").Value = 46
$ws.Range("L25").Value = 104
$ws.Range("L29").Value = 993
$ws.Range("L33").Value = 817
$ws.Range("L34").Value = 104
$ws.Range("L36").Value = 225
$ws.Range("K37").Value = 909
$ws.Range("L37").Value = 669
$ws.Range("L42").Value = 572
$ws.Range("L47").Value = 118
$ws.Range("L49").Value = 88
$ws.Range("L52").Value = 355
$ws.Range("L53").Value = 194
$ws.Range("L54").Value = 384
$ws.Range("L60").Value = 114
$ws.Range("L63").Value = 51
$ws.Range("L64").Value = 117
$ws.Range("L65").Value = 342
$ws.Range("L67").Value = 603
$ws.Range("L73").Value = 140
$ws.Range("L78").Value = 224
$ws.Range("L80").Value = 57
$ws.Range("L83").Value = 389
$ws.Range("L84").Value = 172
$ws.Range("L85").Value = 882
$ws.Range("L90").Value = 182
$ws.Range("L95").Value = 247
$ws.Range("L96").Value = 200
$ws.Range("K101").Value = 27584
$ws.Range("L101").Value = 17640

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L4").Value = 12
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L6").Value = 136
$ws.Range("L7").Value = 574

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 337
$ws.Range("L6").Value = 298

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L6").Value = 71
$ws.Range("L7").Value = 291

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 49
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 136

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L4").Value = 21
$ws.Range("L6").Value = 132
$ws.Range("L7").Value = 479

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 137
$ws.Range("L3").Value = 142
$ws.Range("L7").Value = 436

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L6").Value = 10
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 380
$ws.Range("L7").Value = 993

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 226
$ws.Range("L3").Value = 283
$ws.Range("L7").Value = 817

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L3").Value = 30
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 79
$ws.Range("L5").Value = 2
$ws.Range("L7").Value = 225

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 259
$ws.Range("L2").Value = 201
$ws.Range("K7").Value = 909
$ws.Range("L7").Value = 669

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 195
$ws.Range("L7").Value = 572

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 39
$ws.Range("L7").Value = 118

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L6").Value = 35
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 116
$ws.Range("L6").Value = 94
$ws.Range("L7").Value = 355

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 194

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 186
$ws.Range("L7").Value = 384

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L2").Value = 36
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L2").Value = 34
$ws.Range("L7").Value = 117

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L4").Value = 18
$ws.Range("L7").Value = 342

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 179
$ws.Range("L3").Value = 232
$ws.Range("L7").Value = 603

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 62
$ws.Range("L3").Value = 72
$ws.Range("L4").Value = 22
$ws.Range("L7").Value = 224

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 124
$ws.Range("L4").Value = 15
$ws.Range("L6").Value = 86
$ws.Range("L7").Value = 389

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L6").Value = 48
$ws.Range("L7").Value = 172

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 361
$ws.Range("L6").Value = 183
$ws.Range("L7").Value = 882

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 61
$ws.Range("L7").Value = 182

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L6").Value = 55
$ws.Range("L7").Value = 247

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 200
